$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  0.09375,
  0.1875,
  0.28125,
  0.375,
  0.4687499999999999,
  0.5625,
  0.65625,
  0.75,
  0.84375,
  0.9374999999999999,
  1.03125,
  1.125,
  1.21875,
  1.3125,
  1.40625,
  1.5,
  1.59375,
  1.6875,
  1.78125,
  1.875,
  1.96875,
  2.0625,
  2.15625,
  2.25,
  2.34375,
  2.4375,
  2.53125,
  2.625,
  2.71875,
  2.8125,
  2.90625,
  3,
  3.09375,
  3.1875,
  3.28125,
  3.375,
  3.46875,
  3.5625,
  3.65625,
  3.75,
  3.84375,
  3.9375,
  4.03125,
  4.125,
  4.21875,
  4.3125,
  4.40625,
  4.5,
  4.59375,
  4.6875,
  4.78125,
  4.875,
  4.96875,
  5.0625,
  5.156249999999999,
  5.25,
  5.34375,
  5.4375,
  5.531250000000001,
  5.625,
  5.71875,
  5.812499999999999,
  5.90625
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
